$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:67 down to 44:68.
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new weekly price observation.
$ws.Cells.Item(43, 1).Value = 11
$ws.Cells.Item(43, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(43, 3).Value = "Bíobío"
$ws.Cells.Item(43, 4).Value = 44455
$ws.Cells.Item(43, 5).Value = 8
$ws.Cells.Item(43, 6).Value = 100112043
$ws.Cells.Item(43, 7).Value = "Pepino ensalada"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 100
$ws.Cells.Item(43, 11).Value = 15000
$ws.Cells.Item(43, 12).Value = 16000
$ws.Cells.Item(43, 13).Value = 15500
$ws.Cells.Item(43, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(43, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value = 310
$ws.Cells.Item(43, 17).Value = 50
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Ensure date column keeps the same number format style as the rest of column D.
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
